$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hdi")
$ws.Activate()

# --- Column widths ---
# (ColumnWidth is in "characters"; the engine stores width in a slightly
# larger internal unit, so the COM values below are chosen so the saved
# OOXML <col width="..."/> comes out to exactly 20 / 6.)
$ws.Columns.Item(2).ColumnWidth = 19.17
$ws.Columns.Item(5).ColumnWidth = 5.17

# --- Fix tie-break ordering among countries sharing the same HDI rank ---
# (each block below swaps the B:J values between two rows that tied on column A)
# Row 33
$ws.Cells.Item(33, 2).Value = 'Qatar'
$ws.Cells.Item(33, 3).Value = 25.354825999999999
$ws.Cells.Item(33, 4).Value = 51.183883999999999
$ws.Cells.Item(33, 5).Value = 0.85
$ws.Cells.Item(33, 6).Value = 78.2
$ws.Cells.Item(33, 7).Value = 13.8
$ws.Cells.Item(33, 8).Value = 9.1
$ws.Cells.Item(33, 9).Value = 123124
$ws.Cells.Item(33, 10).Value = -31

# Row 34
$ws.Cells.Item(34, 2).Value = 'Cyprus'
$ws.Cells.Item(34, 3).Value = 35.126412999999999
$ws.Cells.Item(34, 4).Value = 33.429859
$ws.Cells.Item(34, 5).Value = 0.85
$ws.Cells.Item(34, 6).Value = 80.2
$ws.Cells.Item(34, 7).Value = 14
$ws.Cells.Item(34, 8).Value = 11.6
$ws.Cells.Item(34, 9).Value = 28633
$ws.Cells.Item(34, 10).Value = 3

# Row 38
$ws.Cells.Item(38, 2).Value = 'Malta'
$ws.Cells.Item(38, 3).Value = 35.937496000000003
$ws.Cells.Item(38, 4).Value = 14.375416
$ws.Cells.Item(38, 5).Value = 0.83899999999999997
$ws.Cells.Item(38, 6).Value = 80.599999999999994
$ws.Cells.Item(38, 7).Value = 14.4
$ws.Cells.Item(38, 8).Value = 10.3
$ws.Cells.Item(38, 9).Value = 27930
$ws.Cells.Item(38, 10).Value = -1

# Row 39
$ws.Cells.Item(39, 2).Value = 'Lithuania'
$ws.Cells.Item(39, 3).Value = 55.169438
$ws.Cells.Item(39, 4).Value = 23.881274999999999
$ws.Cells.Item(39, 5).Value = 0.83899999999999997
$ws.Cells.Item(39, 6).Value = 73.3
$ws.Cells.Item(39, 7).Value = 16.399999999999999
$ws.Cells.Item(39, 8).Value = 12.4
$ws.Cells.Item(39, 9).Value = 24500
$ws.Cells.Item(39, 10).Value = 7

# Row 51
$ws.Cells.Item(51, 2).Value = 'Russian Federation'
$ws.Cells.Item(51, 3).Value = 61.524009999999997
$ws.Cells.Item(51, 4).Value = 105.31875599999999
$ws.Cells.Item(51, 5).Value = 0.79800000000000004
$ws.Cells.Item(51, 6).Value = 70.099999999999994
$ws.Cells.Item(51, 7).Value = 14.7
$ws.Cells.Item(51, 8).Value = 12
$ws.Cells.Item(51, 9).Value = 22352
$ws.Cells.Item(51, 10).Value = -1

# Row 52
$ws.Cells.Item(52, 2).Value = 'Belarus'
$ws.Cells.Item(52, 3).Value = 53.709806999999998
$ws.Cells.Item(52, 4).Value = 27.953389000000001
$ws.Cells.Item(52, 5).Value = 0.79800000000000004
$ws.Cells.Item(52, 6).Value = 71.3
$ws.Cells.Item(52, 7).Value = 15.7
$ws.Cells.Item(52, 8).Value = 12
$ws.Cells.Item(52, 9).Value = 16676
$ws.Cells.Item(52, 10).Value = 14

# Row 54
$ws.Cells.Item(54, 2).Value = 'Uruguay'
$ws.Cells.Item(54, 3).Value = -32.522779
$ws.Cells.Item(54, 4).Value = -55.765835000000003
$ws.Cells.Item(54, 5).Value = 0.79300000000000004
$ws.Cells.Item(54, 6).Value = 77.2
$ws.Cells.Item(54, 7).Value = 15.5
$ws.Cells.Item(54, 8).Value = 8.5
$ws.Cells.Item(54, 9).Value = 19283
$ws.Cells.Item(54, 10).Value = 7

# Row 55
$ws.Cells.Item(55, 2).Value = 'Romania'
$ws.Cells.Item(55, 3).Value = 45.943161000000003
$ws.Cells.Item(55, 4).Value = 24.966760000000001
$ws.Cells.Item(55, 5).Value = 0.79300000000000004
$ws.Cells.Item(55, 6).Value = 74.7
$ws.Cells.Item(55, 7).Value = 14.2
$ws.Cells.Item(55, 8).Value = 10.8
$ws.Cells.Item(55, 9).Value = 18108
$ws.Cells.Item(55, 10).Value = 10

# Row 61
$ws.Cells.Item(61, 2).Value = 'Panama'
$ws.Cells.Item(61, 3).Value = 8.5379810000000003
$ws.Cells.Item(61, 4).Value = -80.782127000000003
$ws.Cells.Item(61, 5).Value = 0.78
$ws.Cells.Item(61, 6).Value = 77.599999999999994
$ws.Cells.Item(61, 7).Value = 13.3
$ws.Cells.Item(61, 8).Value = 9.3000000000000007
$ws.Cells.Item(61, 9).Value = 18192
$ws.Cells.Item(61, 10).Value = 1

# Row 62
$ws.Cells.Item(62, 2).Value = 'Palau'
$ws.Cells.Item(62, 3).Value = 7.5149800000000004
$ws.Cells.Item(62, 4).Value = 134.58251999999999
$ws.Cells.Item(62, 5).Value = 0.78
$ws.Cells.Item(62, 6).Value = 72.7
$ws.Cells.Item(62, 7).Value = 13.7
$ws.Cells.Item(62, 8).Value = 12.3
$ws.Cells.Item(62, 9).Value = 13496
$ws.Cells.Item(62, 10).Value = 18

# Row 65
$ws.Cells.Item(65, 2).Value = 'Trinidad and Tobago'
$ws.Cells.Item(65, 3).Value = 10.691803
$ws.Cells.Item(65, 4).Value = -61.222503000000003
$ws.Cells.Item(65, 5).Value = 0.77200000000000002
$ws.Cells.Item(65, 6).Value = 70.400000000000006
$ws.Cells.Item(65, 7).Value = 12.3
$ws.Cells.Item(65, 8).Value = 10.9
$ws.Cells.Item(65, 9).Value = 26090
$ws.Cells.Item(65, 10).Value = -25

# Row 66
$ws.Cells.Item(66, 2).Value = 'Seychelles'
$ws.Cells.Item(66, 3).Value = -4.6795739999999997
$ws.Cells.Item(66, 4).Value = 55.491976999999999
$ws.Cells.Item(66, 5).Value = 0.77200000000000002
$ws.Cells.Item(66, 6).Value = 73.099999999999994
$ws.Cells.Item(66, 7).Value = 13.4
$ws.Cells.Item(66, 8).Value = 9.4
$ws.Cells.Item(66, 9).Value = 23300
$ws.Cells.Item(66, 10).Value = -19

# Row 68
$ws.Cells.Item(68, 2).Value = 'Lebanon'
$ws.Cells.Item(68, 3).Value = 33.854720999999998
$ws.Cells.Item(68, 4).Value = 35.862285
$ws.Cells.Item(68, 5).Value = 0.76900000000000002
$ws.Cells.Item(68, 6).Value = 79.3
$ws.Cells.Item(68, 7).Value = 13.8
$ws.Cells.Item(68, 8).Value = 7.9
$ws.Cells.Item(68, 9).Value = 16509
$ws.Cells.Item(68, 10).Value = -1

# Row 69
$ws.Cells.Item(69, 2).Value = 'Cuba'
$ws.Cells.Item(69, 3).Value = 21.521757000000001
$ws.Cells.Item(69, 4).Value = -77.781166999999996
$ws.Cells.Item(69, 5).Value = 0.76900000000000002
$ws.Cells.Item(69, 6).Value = 79.400000000000006
$ws.Cells.Item(69, 7).Value = 13.8
$ws.Cells.Item(69, 8).Value = 11.5
$ws.Cells.Item(69, 9).Value = 7301
$ws.Cells.Item(69, 10).Value = 47

# Row 70
$ws.Cells.Item(70, 2).Value = 'Iran (Islamic Republic of)'
$ws.Cells.Item(70, 3).Value = 32.427908000000002
$ws.Cells.Item(70, 4).Value = 53.688046
$ws.Cells.Item(70, 5).Value = 0.76600000000000001
$ws.Cells.Item(70, 6).Value = 75.400000000000006
$ws.Cells.Item(70, 7).Value = 15.1
$ws.Cells.Item(70, 8).Value = 8.1999999999999993
$ws.Cells.Item(70, 9).Value = 15440
$ws.Cells.Item(70, 10).Value = 4

# Row 71
$ws.Cells.Item(71, 2).Value = 'Costa Rica'
$ws.Cells.Item(71, 3).Value = 9.7489170000000005
$ws.Cells.Item(71, 4).Value = -83.753428
$ws.Cells.Item(71, 5).Value = 0.76600000000000001
$ws.Cells.Item(71, 6).Value = 79.400000000000006
$ws.Cells.Item(71, 7).Value = 13.9
$ws.Cells.Item(71, 8).Value = 8.4
$ws.Cells.Item(71, 9).Value = 13413
$ws.Cells.Item(71, 10).Value = 10

# Row 87
$ws.Cells.Item(87, 2).Value = 'Bosnia and Herzegovina'
$ws.Cells.Item(87, 3).Value = 43.915886
$ws.Cells.Item(87, 4).Value = 17.679075999999998
$ws.Cells.Item(87, 5).Value = 0.73299999999999998
$ws.Cells.Item(87, 6).Value = 76.5
$ws.Cells.Item(87, 7).Value = 13.6
$ws.Cells.Item(87, 8).Value = 8.3000000000000007
$ws.Cells.Item(87, 9).Value = 9638
$ws.Cells.Item(87, 10).Value = 19

# Row 88
$ws.Cells.Item(88, 2).Value = 'Armenia'
$ws.Cells.Item(88, 3).Value = 40.069099000000001
$ws.Cells.Item(88, 4).Value = 45.038189000000003
$ws.Cells.Item(88, 5).Value = 0.73299999999999998
$ws.Cells.Item(88, 6).Value = 74.7
$ws.Cells.Item(88, 7).Value = 12.3
$ws.Cells.Item(88, 8).Value = 10.9
$ws.Cells.Item(88, 9).Value = 8124
$ws.Cells.Item(88, 10).Value = 22

# Row 92
$ws.Cells.Item(92, 2).Value = 'Mongolia'
$ws.Cells.Item(92, 3).Value = 46.862496
$ws.Cells.Item(92, 4).Value = 103.846656
$ws.Cells.Item(92, 5).Value = 0.72699999999999998
$ws.Cells.Item(92, 6).Value = 69.400000000000006
$ws.Cells.Item(92, 7).Value = 14.6
$ws.Cells.Item(92, 8).Value = 9.3000000000000007
$ws.Cells.Item(92, 9).Value = 10729
$ws.Cells.Item(92, 10).Value = 4

# Row 93
$ws.Cells.Item(93, 2).Value = 'Fiji'
$ws.Cells.Item(93, 3).Value = -16.578192999999999
$ws.Cells.Item(93, 4).Value = 179.414413
$ws.Cells.Item(93, 5).Value = 0.72699999999999998
$ws.Cells.Item(93, 6).Value = 70
$ws.Cells.Item(93, 7).Value = 15.7
$ws.Cells.Item(93, 8).Value = 9.9
$ws.Cells.Item(93, 9).Value = 7493
$ws.Cells.Item(93, 10).Value = 21

# Row 95
$ws.Cells.Item(95, 2).Value = 'Libya'
$ws.Cells.Item(95, 3).Value = 26.335100000000001
$ws.Cells.Item(95, 4).Value = 17.228331000000001
$ws.Cells.Item(95, 5).Value = 0.72399999999999998
$ws.Cells.Item(95, 6).Value = 71.599999999999994
$ws.Cells.Item(95, 7).Value = 14
$ws.Cells.Item(95, 8).Value = 7.3
$ws.Cells.Item(95, 9).Value = 14911
$ws.Cells.Item(95, 10).Value = -19

# Row 96
$ws.Cells.Item(96, 2).Value = 'Dominica'
$ws.Cells.Item(96, 3).Value = 15.414999
$ws.Cells.Item(96, 4).Value = -61.370975999999999
$ws.Cells.Item(96, 5).Value = 0.72399999999999998
$ws.Cells.Item(96, 6).Value = 77.8
$ws.Cells.Item(96, 7).Value = 12.7
$ws.Cells.Item(96, 8).Value = 7.9
$ws.Cells.Item(96, 9).Value = 9994
$ws.Cells.Item(96, 10).Value = 4

# Row 102
$ws.Cells.Item(102, 2).Value = 'Dominican Republic'
$ws.Cells.Item(102, 3).Value = 18.735693000000001
$ws.Cells.Item(102, 4).Value = -70.162650999999997
$ws.Cells.Item(102, 5).Value = 0.71499999999999997
$ws.Cells.Item(102, 6).Value = 73.5
$ws.Cells.Item(102, 7).Value = 13.1
$ws.Cells.Item(102, 8).Value = 7.6
$ws.Cells.Item(102, 9).Value = 11883
$ws.Cells.Item(102, 10).Value = -12

# Row 103
$ws.Cells.Item(103, 2).Value = 'Belize'
$ws.Cells.Item(103, 3).Value = 17.189876999999999
$ws.Cells.Item(103, 4).Value = -88.497649999999993
$ws.Cells.Item(103, 5).Value = 0.71499999999999997
$ws.Cells.Item(103, 6).Value = 70
$ws.Cells.Item(103, 7).Value = 13.6
$ws.Cells.Item(103, 8).Value = 10.5
$ws.Cells.Item(103, 9).Value = 7614
$ws.Cells.Item(103, 10).Value = 9

# Row 117
$ws.Cells.Item(117, 2).Value = 'South Africa'
$ws.Cells.Item(117, 3).Value = -30.559481999999999
$ws.Cells.Item(117, 4).Value = 22.937505999999999
$ws.Cells.Item(117, 5).Value = 0.66600000000000004
$ws.Cells.Item(117, 6).Value = 57.4
$ws.Cells.Item(117, 7).Value = 13.6
$ws.Cells.Item(117, 8).Value = 9.9
$ws.Cells.Item(117, 9).Value = 12122
$ws.Cells.Item(117, 10).Value = -29

# Row 118
$ws.Cells.Item(118, 2).Value = 'El Salvador'
$ws.Cells.Item(118, 3).Value = 13.794185000000001
$ws.Cells.Item(118, 4).Value = -88.896529999999998
$ws.Cells.Item(118, 5).Value = 0.66600000000000004
$ws.Cells.Item(118, 6).Value = 73
$ws.Cells.Item(118, 7).Value = 12.3
$ws.Cells.Item(118, 8).Value = 6.5
$ws.Cells.Item(118, 9).Value = 7349
$ws.Cells.Item(118, 10).Value = -3

# Row 127
$ws.Cells.Item(127, 2).Value = 'Namibia'
$ws.Cells.Item(127, 3).Value = -22.957640000000001
$ws.Cells.Item(127, 4).Value = 18.490410000000001
$ws.Cells.Item(127, 5).Value = 0.628
$ws.Cells.Item(127, 6).Value = 64.8
$ws.Cells.Item(127, 7).Value = 11.3
$ws.Cells.Item(127, 8).Value = 6.2
$ws.Cells.Item(127, 9).Value = 9418
$ws.Cells.Item(127, 10).Value = -21

# Row 128
$ws.Cells.Item(128, 2).Value = 'Morocco'
$ws.Cells.Item(128, 3).Value = 31.791702000000001
$ws.Cells.Item(128, 4).Value = -7.0926200000000001
$ws.Cells.Item(128, 5).Value = 0.628
$ws.Cells.Item(128, 6).Value = 74
$ws.Cells.Item(128, 7).Value = 11.6
$ws.Cells.Item(128, 8).Value = 4.4000000000000004
$ws.Cells.Item(128, 9).Value = 6850
$ws.Cells.Item(128, 10).Value = -8

# Row 135
$ws.Cells.Item(135, 2).Value = 'Vanuatu'
$ws.Cells.Item(135, 3).Value = -15.376706
$ws.Cells.Item(135, 4).Value = 166.959158
$ws.Cells.Item(135, 5).Value = 0.59399999999999997
$ws.Cells.Item(135, 6).Value = 71.900000000000006
$ws.Cells.Item(135, 7).Value = 10.6
$ws.Cells.Item(135, 8).Value = 6.8
$ws.Cells.Item(135, 9).Value = 2803
$ws.Cells.Item(135, 10).Value = 19

# Row 136
$ws.Cells.Item(136, 2).Value = 'Syrian Arab Republic'
$ws.Cells.Item(136, 3).Value = 34.802075000000002
$ws.Cells.Item(136, 4).Value = 38.996814999999998
$ws.Cells.Item(136, 5).Value = 0.59399999999999997
$ws.Cells.Item(136, 6).Value = 69.599999999999994
$ws.Cells.Item(136, 7).Value = 12.3
$ws.Cells.Item(136, 8).Value = 6.3
$ws.Cells.Item(136, 9).Value = 2728
$ws.Cells.Item(136, 10).Value = 21

# Row 165
$ws.Cells.Item(165, 2).Value = 'Uganda'
$ws.Cells.Item(165, 3).Value = 1.3733329999999999
$ws.Cells.Item(165, 4).Value = 32.290275000000001
$ws.Cells.Item(165, 5).Value = 0.48299999999999998
$ws.Cells.Item(165, 6).Value = 58.5
$ws.Cells.Item(165, 7).Value = 9.8000000000000007
$ws.Cells.Item(165, 8).Value = 5.4
$ws.Cells.Item(165, 9).Value = 1613
$ws.Cells.Item(165, 10).Value = 6

# Row 166
$ws.Cells.Item(166, 2).Value = 'Rwanda'
$ws.Cells.Item(166, 3).Value = -1.9402779999999999
$ws.Cells.Item(166, 4).Value = 29.873888000000001
$ws.Cells.Item(166, 5).Value = 0.48299999999999998
$ws.Cells.Item(166, 6).Value = 64.2
$ws.Cells.Item(166, 7).Value = 10.3
$ws.Cells.Item(166, 8).Value = 9.8000000000000007
$ws.Cells.Item(166, 9).Value = 1458
$ws.Cells.Item(166, 10).Value = 11

# --- Re-apply sort on column A (keeps sortState bookkeeping in sync with the data) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A189"))
$ws.Sort.SetRange($ws.Range("A1:J189"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Selection state ---
$ws.Range("K13").Select()
